$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (original row 26) entirely - all rows below shift up by one.
$ws.Rows(26).Delete()

# Delete the "SC 92" row (original row 28, now row 27 after the first delete) entirely.
$ws.Rows(27).Delete()

# Clear / set individual data cells (values imputed/removed) to match the target state.
$ws.Range("E2").ClearContents()
$ws.Range("D6").Value = -14.2
$ws.Range("D8").ClearContents()
$ws.Range("D18").Value = -15.2
$ws.Range("D20").ClearContents()
$ws.Range("D23").Value = -13.9
$ws.Range("D25").ClearContents()

$ws.Range("B27").Value = -20.4
$ws.Range("B28").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B30").Value = -19.7
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("B32").ClearContents()
